# Update the "runs/balls/fours/sixes" stats for rows 2-13 on the
# "Aaron Finch " sheet. Values in these columns are stored as text
# (numbers-stored-as-text), so force a text number format before writing
# the values to keep their string type intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("13", "14", "1", "0")
    3  = @("20", "18", "2", "1")
    4  = @("8",  "7",  "2", "0")
    5  = @("15", "11", "3", "0")
    6  = @("14", "11", "0", "2")
    7  = @("52", "35", "7", "1")
    8  = @("16", "21", "2", "0")
    9  = @("32", "30", "3", "1")
    10 = @("47", "37", "4", "1")
    11 = @("2",  "9",  "0", "0")
    12 = @("29", "27", "1", "2")
    13 = @("20", "21", "3", "0")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]

    $rng = $ws.Range("C$row`:F$row")
    $rng.NumberFormat = "@"

    $ws.Range("C$row").Value = $values[0]
    $ws.Range("D$row").Value = $values[1]
    $ws.Range("E$row").Value = $values[2]
    $ws.Range("F$row").Value = $values[3]
}
